$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Ginnie@123"
$ws.Range("B1").Select()
